$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.729336738586426
$ws.Range("B1").Value = 5.812500953674316
$ws.Range("C1").Value = 5.600942134857178
$ws.Range("D1").Value = 9.303844451904297
$ws.Range("E1").Value = 7.044075012207031
